{"js": "// Apply README/docx stats fix-up for the Renaissance / JDK21 / ZGC\n// neo4j-analytics 1G benchmark table. The document's single table has\n// one value per row; update the affected rows in place.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// rowIndex -> new cell text\nconst updates = {\n  0: \"0M\",\n  1: \"0M\",\n  2: \"0M\",\n  3: \"75\",\n  6: \"0.30899\",\n  7: \"0.06737\",\n  11: \"15.87958\",\n  43: \"6.7\",\n  44: \"15.88\",\n  45: \"17\",\n};\n\nfor (const rowIndex of Object.keys(updates)) {\n  const idx = parseInt(rowIndex, 10);\n  const cell = table.getCell(idx, 0);\n  cell.value = updates[rowIndex];\n}\n\nawait context.sync();\n", "ps1": "# Apply README/docx stats fix-up for the Renaissance / JDK21 / ZGC\n# neo4j-analytics 1G benchmark table. The document's single table has\n# one value per row; update the affected rows in place.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Word COM tables/cells are 1-based.\n$updates = @{\n    1  = \"0M\"\n    2  = \"0M\"\n    3  = \"0M\"\n    4  = \"75\"\n    7  = \"0.30899\"\n    8  = \"0.06737\"\n    12 = \"15.87958\"\n    44 = \"6.7\"\n    45 = \"15.88\"\n    46 = \"17\"\n}\n\nforeach ($rowNum in $updates.Keys) {\n    $cell = $t.Cell($rowNum, 1)\n    $cell.Range.Text = $updates[$rowNum]\n}\n"}
